# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`r`n✅ Dólar paralelo: 68`r`n`r`nBinance`r`n✅ 1000 Bs = 13.12 = 53005.25 pesos`r`n✅ 53005.25 pesos = 13.06 = 974.62 Bs`r`n`r`nPromedio competencia`r`n✅ Tasa pesos: 20`r`n✅ Tasa Bs: 20`r`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- tasas: update the rate values in N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 76.2
$ws2.Range("O10").Value = 4039
$ws2.Range("N12").Value = 4057.99
$ws2.Range("O12").Value = 74.615
